# Re-apply "Set precision as displayed" style rounding to the Longitude/
# Latitude columns (H, I) for rows 5-19, add a few previously-blank flag
# cells (N, O, Q columns) that were filled in, reset row heights that had
# been manually set back to automatic (matches rows whose wrapped
# description text no longer needs the extra height), and restore the
# sheet scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Toiletten")

# --- Row 5 ---------------------------------------------------------------
$ws.Range("H5").Value = 13.388736093
$ws.Range("I5").Value = 52.548913458000001
$ws.Range("N5").Value = 0
$ws.Range("Q5").Value = 6

# --- Row 6 ---------------------------------------------------------------
$ws.Range("H6").Value = 13.38851667
$ws.Range("I6").Value = 52.549208999999998
$ws.Range("N6").Value = 0
$ws.Range("Q6").Value = 6

# --- Row 7 ---------------------------------------------------------------
$ws.Range("H7").Value = 13.334153323000001
$ws.Range("I7").Value = 52.508296522999999
$ws.Range("N7").Value = 0
$ws.Range("Q7").Value = 6
$ws.Rows.Item(7).AutoFit()

# --- Row 8 ---------------------------------------------------------------
$ws.Range("H8").Value = 13.387249924000001
$ws.Range("I8").Value = 52.519966529999998
$ws.Range("N8").Value = 0
$ws.Range("Q8").Value = 6
$ws.Rows.Item(8).AutoFit()

# --- Row 9 ---------------------------------------------------------------
$ws.Range("H9").Value = 13.411754707
$ws.Range("I9").Value = 52.521398576000003
$ws.Range("N9").Value = 0
$ws.Range("Q9").Value = 6
$ws.Rows.Item(9).AutoFit()

# --- Row 10 --------------------------------------------------------------
$ws.Range("H10").Value = 13.191026179
$ws.Range("I10").Value = 52.434341541000002
$ws.Range("N10").Value = 0
$ws.Range("Q10").Value = 6
$ws.Rows.Item(10).AutoFit()

# --- Row 11 --------------------------------------------------------------
$ws.Range("H11").Value = 13.467933001
$ws.Range("I11").Value = 52.503972644000001
$ws.Range("N11").Value = 0
$ws.Range("Q11").Value = 6
$ws.Rows.Item(11).AutoFit()

# --- Row 12 (previously had no M/N/P/Q/R values at all) ------------------
$ws.Range("H12").Value = 13.497256016
$ws.Range("I12").Value = 52.509409104
$ws.Range("M12").Value = 1
$ws.Range("N12").Value = 0
$ws.Range("P12").Value = 1
$ws.Range("Q12").Value = 6
$ws.Range("R12").Value = 1
$ws.Rows.Item(12).AutoFit()

# --- Row 13 --------------------------------------------------------------
$ws.Range("H13").Value = 13.179296329
$ws.Range("I13").Value = 52.421585714999999
$ws.Range("N13").Value = 0
$ws.Range("Q13").Value = 6
$ws.Rows.Item(13).AutoFit()

# --- Row 14 (row height untouched) ---------------------------------------
$ws.Range("H14").Value = 13.199195979000001
$ws.Range("I14").Value = 52.534650505999998
$ws.Range("N14").Value = 0
$ws.Range("Q14").Value = 6

# --- Row 15 (row height untouched) ---------------------------------------
$ws.Range("H15").Value = 13.365804841999999
$ws.Range("I15").Value = 52.476169622
$ws.Range("N15").Value = 0
$ws.Range("Q15").Value = 6

# --- Row 16 (row height untouched) ---------------------------------------
$ws.Range("H16").Value = 13.364173495999999
$ws.Range("I16").Value = 52.476434922999999
$ws.Range("N16").Value = 0
$ws.Range("Q16").Value = 6

# --- Row 17 (row height untouched) ---------------------------------------
$ws.Range("H17").Value = 13.435104375
$ws.Range("I17").Value = 52.510225527999999
$ws.Range("N17").Value = 0
$ws.Range("Q17").Value = 6

# --- Row 18 ----------------------------------------------------------------
$ws.Range("H18").Value = 13.375946813000001
$ws.Range("I18").Value = 52.509646781000001
$ws.Range("N18").Value = 0
$ws.Range("Q18").Value = 6
$ws.Rows.Item(18).AutoFit()

# --- Row 19 (row height untouched, also gets O19) -------------------------
$ws.Range("H19").Value = 13.368468333999999
$ws.Range("I19").Value = 52.525689796999998
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 1
$ws.Range("Q19").Value = 6

# --- View state: scroll back to the top-left corner and select B21:B22 ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B21:B22").Select() | Out-Null
